$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44/45: swap EnergySwap and Frax (with updated link/values) ---
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'

# --- Force text format on Price cells whose new values would otherwise be
#     auto-detected as numbers by Excel (losing the exact text formatting) ---
$textCells = @("D4","D7","D8","D9","D10","D11","D13","D14","D15","D17","D18","D19","D21","D22","D23","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

# --- Update Price (D) and Volume 1h (E) columns ---
$ws.Range("D2").Value = '22.465.99'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '1.572.17'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.3752'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '49.88'
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").Value = '0.3410'
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D10").Value = '1.150'
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("D11").Value = '0.07566'
$ws.Range("E11").Value = '  -1.52%  '
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '21.39'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").Value = '6.019'
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").Value = '6.960'
$ws.Range("E15").Value = '  +0.40%  '
$ws.Range("D16").Value = '1.575.22'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '0.00001122'
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").Value = '90.89'
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").Value = '0.06746'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '6.263'
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").Value = '16.41'
$ws.Range("E22").Value = '  -2.29%  '
$ws.Range("D23").Value = '12.23'
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("D24").Value = '22.455.59'
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("D25").Value = '2.355'
$ws.Range("E25").Value = '  -1.56%  '
$ws.Range("D26").Value = '2.600'
$ws.Range("E26").Value = '  -6.36%  '
$ws.Range("D27").Value = '20.15'
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").Value = '148.59'
$ws.Range("E28").Value = '  +2.81%  '
$ws.Range("D29").Value = '5.003'
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("D30").Value = '126.04'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").Value = '1.750.33'
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").Value = '1.034'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("D33").Value = '6.145'
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").Value = '1.986'
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("D35").Value = '9.863'
$ws.Range("E35").Value = '  -2.08%  '
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").Value = '1.377'
$ws.Range("E37").Value = '  +6.45%  '
$ws.Range("D38").Value = '0.02467'
$ws.Range("E38").Value = '  -3.65%  '
$ws.Range("D39").Value = '0.2293'
$ws.Range("E39").Value = '  -1.42%  '
$ws.Range("D40").Value = '0.06587'
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").Value = '5.482'
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("D42").Value = '11.37'
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("D43").Value = '0.6309'
$ws.Range("E43").Value = '  -2.13%  '
$ws.Range("D44").Value = '14.10'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '3.817'
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("D47").Value = '0.5880'
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("D48").Value = '2.104'
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").Value = '130.08'
$ws.Range("E49").Value = '  +4.04%  '
$ws.Range("D50").Value = '1.226'
$ws.Range("E50").Value = '  -5.86%  '
$ws.Range("D51").Value = '0.07327'
$ws.Range("E51").Value = '  -0.14%  '
